$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new "Price" values look like plain decimal numbers (e.g. "488.60").
# Assigning such a string straight to .Value makes Excel auto-convert it into a numeric
# cell, which would lose the exact original text (e.g. trailing zero, or flip it into a
# floating point number). To guarantee the literal text is stored (matching the original
# inline-string cells), we build it as a quoted-text formula ( ="488.60" ) and then convert
# that formula result to a plain value in place via copy / paste-special-values. This keeps
# the cell a plain text cell without touching any number formats or styles.
function Set-CellText($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $escaped = $val -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
}

$ws.Range("D2").Value = "68.918.73"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "3.930.47"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-CellText "D5" "488.60"
$ws.Range("E5").Value = "  +0.44%  "
Set-CellText "D6" "146.47"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("E11").Value = "  -4.74%  "
Set-CellText "D12" "42.96"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").Value = "4.557.45"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "3.934.29"
$ws.Range("E15").Value = "  -0.11%  "
Set-CellText "D16" "14.23"
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("E17").Value = "  -0.66%  "
Set-CellText "D18" "19.92"
$ws.Range("E18").Value = "  +0.04%  "
Set-CellText "D19" "1.16"
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("D20").Value = "69.008.36"
$ws.Range("E20").Value = "  +1.02%  "
Set-CellText "D21" "436.08"
$ws.Range("E21").Value = "  -2.28%  "
Set-CellText "D22" "3.46"
$ws.Range("E22").Value = "  +2.41%  "
Set-CellText "D23" "14.52"
$ws.Range("E23").Value = "  -2.19%  "
Set-CellText "D24" "12.45"
$ws.Range("E24").Value = "  +14.77%  "
Set-CellText "D25" "89.39"
$ws.Range("E25").Value = "  +0.73%  "
Set-CellText "D26" "3.73"
$ws.Range("E26").Value = "  +3.23%  "
Set-CellText "D27" "11.11"
$ws.Range("E27").Value = "  -2.71%  "
Set-CellText "D28" "37.14"
$ws.Range("E28").Value = "  -4.22%  "
Set-CellText "D29" "5.66"
$ws.Range("E29").Value = "  -3.75%  "
Set-CellText "D30" "709.62"
$ws.Range("E30").Value = "  +2.68%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-CellText "D31" "0.133"
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-CellText "D32" "13.51"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  +2.93%  "
Set-CellText "D34" "0.479"
$ws.Range("E34").Value = "  +29.76%  "
$ws.Range("D35").Value = "0.0₃0892"
$ws.Range("E35").Value = "  -5.93%  "
Set-CellText "D36" "61.82"
$ws.Range("E36").Value = "  +4.60%  "
$ws.Range("E37").Value = "  +6.20%  "
Set-CellText "D38" "40.76"
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("E39").Value = "  -0.19%  "
Set-CellText "D40" "1.00"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("E43").Value = "  +2.82%  "
Set-CellText "D44" "3.06"
$ws.Range("E44").Value = "  -2.08%  "
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("E47").Value = "  +5.59%  "
$ws.Range("D48").Value = "0.0₆0362"
$ws.Range("E48").Value = "  +5.76%  "
Set-CellText "D49" "2.99"
$ws.Range("E49").Value = "  +5.09%  "
Set-CellText "D50" "3.38"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("E51").Value = "  -3.19%  "

$excel.CutCopyMode = 0
